$d = $word.ActiveDocument

# Locate the paragraph that contains the "We looking for..." submission
# requirements text. It is a single paragraph built from several runs,
# each beginning with a line break (vertical-tab, chr(11)).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*We looking for the following as part of your submission:*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $paraStart = $target.Range.Start
    $paraEnd = $target.Range.End

    # The paragraph mark itself sits in the final 1-character slot of the
    # paragraph range (End-1 .. End). Locate the first line break (kept)
    # and delete everything from just after it through just before the
    # paragraph mark.
    $fullRange = $d.Range($paraStart, $paraEnd)
    $text = $fullRange.Text

    $firstBreak = $text.IndexOf([char]11)
    $deleteStart = $paraStart + $firstBreak + 1
    $deleteEnd = $paraEnd - 1

    if ($deleteStart -lt $deleteEnd) {
        $toDelete = $d.Range($deleteStart, $deleteEnd)
        $toDelete.Delete()
    }
}
